$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 40: OUT_FOLDER_MODELS -> data/classifiers
$ws.Range("A40").Value = "OUT_FOLDER_MODELS"
$ws.Range("B40").Value = "data/classifiers"

# New row 41: a blank, underlined separator cell (adds the underline font/style)
$ws.Range("A41").Font.Underline = $true

# Update selection to match the new active cell
$ws.Range("B40").Select() | Out-Null
